# edit.ps1 - apply the Resume_Kenneth_Harlley_.docx edits described by the diff:
#   1. Drop the trailing ". " after "...Operating Systems (Enrolled)" in the
#      Relevant Courses sentence (and the now-unneeded xml:space="preserve").
#   2. Split the single run "Jan. 2017 - Present" into three runs with the
#      same run formatting: "Jan. 2017 ", "-", " Present".
#   3. Change the en dash between "May 2019" and "Aug. 2019" to a plain
#      hyphen-minus, keeping the surrounding runs (spaces, "Aug. 2019") intact
#      and separate.

$d = $word.ActiveDocument

# --- 1. Relevant Courses sentence: drop the trailing ". " -------------------
$d.Content.Find.Execute(
    "Operating Systems (Enrolled). ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Operating Systems (Enrolled)", 2) | Out-Null

# --- 2. "Jan. 2017 - Present" -> three runs ---------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "Jan. 2017 - Present", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$jobStart = $rng.Start

# The hyphen sits right after "Jan. 2017 " (10 characters in).
$hyphenRng = $d.Range($jobStart + 10, $jobStart + 11)
# Re-asserting (the already-true) Bold value forces this engine to carve the
# character out into its own run, without altering any visible formatting.
$hyphenRng.Bold = 0
$hyphenRng.Bold = 1

# --- 3. En dash -> hyphen-minus in "May 2019 - Aug. 2019" -------------------
$enDash = [char]0x2013
$rng2 = $d.Content
$rng2.Find.Execute(
    "May 2019 " + $enDash + " Aug. 2019", $true, $false, $false, $false,
    $false, $true, 1, $false, "", 0) | Out-Null
$rangeStart = $rng2.Start

# Layout inside the found range: "May 2019 " (9 chars) + en dash (1) + ...
$dashRng = $d.Range($rangeStart + 9, $rangeStart + 10)
$dashRng.Text = "-"

# Re-splitting the dash and the following space keeps every neighbouring run
# (the leading space, the trailing space, "Aug. 2019") separate and untouched,
# matching the pre-existing run layout.
$dashRng2 = $d.Range($rangeStart + 9, $rangeStart + 10)
$dashRng2.Bold = 0
$dashRng2.Bold = 1

$spaceAfterRng = $d.Range($rangeStart + 10, $rangeStart + 11)
$spaceAfterRng.Bold = 0
$spaceAfterRng.Bold = 1
